# cryptos.xlsx symbol-list refresh (GitHub Actions scheduled update)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Hora" column: every data row moves from hour 4 to hour 5
$ws.Range("G2:G51").NumberFormat = "@"
$ws.Range("G2:G51").Value = "5"

# Per-row refreshed figures: Coin, Link, Price, Volume(1h)
# Columns left blank ($null) are unchanged for that row.
$rows = @(
    @{ Row = 2; B = $null; C = $null; D = "309.03"; E = "2.07%" }
    @{ Row = 3; B = $null; C = $null; D = "38.87"; E = "8.90%" }
    @{ Row = 4; B = $null; C = $null; D = "5.094"; E = "1.21%" }
    @{ Row = 5; B = $null; C = $null; D = "0.08185"; E = "3.65%" }
    @{ Row = 6; B = $null; C = $null; D = "2.014"; E = "9.03%" }
    @{ Row = 7; B = "KuCoinToken"; C = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"; D = "7.906"; E = "1.64%" }
    @{ Row = 8; B = "MXToken"; C = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"; D = "0.9312"; E = "1.47%" }
    @{ Row = 9; B = "LiechtensteinCryptoassetsExchange"; C = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"; D = "0.1403"; E = "4.36%" }
    @{ Row = 10; B = "WazirX"; C = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"; D = "0.1949"; E = "3.26%" }
    @{ Row = 11; B = "MandalaExchangeToken"; C = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"; D = "0.09247"; E = "2.29%" }
    @{ Row = 12; B = "BitrueCoin"; C = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"; D = "0.03468"; E = "-0.06%" }
    @{ Row = 13; B = "BitMartToken"; C = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"; D = "0.09850"; E = "0.46%" }
    @{ Row = 14; B = "BitForexToken"; C = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"; D = "0.001422"; E = "1.19%" }
    @{ Row = 15; B = "TigerCash"; C = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"; D = "0.006053"; E = "-1.10%" }
    @{ Row = 16; B = "LEO"; C = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"; D = "3.674"; E = "-1.31%" }
    @{ Row = 17; B = "GateToken"; C = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"; D = "4.178"; E = "1.80%" }
    @{ Row = 18; B = $null; C = $null; D = "3.470"; E = "5.06%" }
    @{ Row = 19; B = $null; C = $null; D = $null; E = "0.41%" }
    @{ Row = 20; B = $null; C = $null; D = $null; E = "-0.53%" }
    @{ Row = 21; B = $null; C = $null; D = "4.804"; E = "-7.39%" }
    @{ Row = 22; B = $null; C = $null; D = $null; E = "11.85%" }
    @{ Row = 23; B = $null; C = $null; D = "0.04459"; E = "1.34%" }
    @{ Row = 24; B = $null; C = $null; D = "0.001244"; E = "0.56%" }
    @{ Row = 25; B = $null; C = $null; D = $null; E = "-9.37%" }
    @{ Row = 26; B = $null; C = $null; D = $null; E = $null }
    @{ Row = 27; B = $null; C = $null; D = $null; E = "-0.03%" }
    @{ Row = 28; B = $null; C = $null; D = $null; E = $null }
    @{ Row = 29; B = $null; C = $null; D = $null; E = $null }
    @{ Row = 30; B = $null; C = $null; D = $null; E = $null }
    @{ Row = 31; B = $null; C = $null; D = $null; E = $null }
    @{ Row = 32; B = $null; C = $null; D = $null; E = $null }
    @{ Row = 33; B = $null; C = $null; D = $null; E = $null }
    @{ Row = 34; B = $null; C = $null; D = $null; E = $null }
    @{ Row = 35; B = $null; C = $null; D = $null; E = $null }
    @{ Row = 36; B = $null; C = $null; D = $null; E = $null }
    @{ Row = 37; B = $null; C = $null; D = $null; E = $null }
    @{ Row = 38; B = $null; C = $null; D = $null; E = $null }
    @{ Row = 39; B = $null; C = $null; D = "0.02133"; E = "10.48%" }
    @{ Row = 40; B = $null; C = $null; D = "0.05171"; E = "-0.47%" }
    @{ Row = 41; B = $null; C = $null; D = "0.007482"; E = "-1.46%" }
    @{ Row = 42; B = $null; C = $null; D = "0.01014"; E = "-0.42%" }
    @{ Row = 43; B = $null; C = $null; D = "0.1368"; E = "2.12%" }
    @{ Row = 44; B = $null; C = $null; D = $null; E = $null }
    @{ Row = 45; B = $null; C = $null; D = "0.009689"; E = "-0.77%" }
    @{ Row = 46; B = $null; C = $null; D = "0.00006318"; E = "2.76%" }
    @{ Row = 47; B = $null; C = $null; D = $null; E = "-0.06%" }
    @{ Row = 48; B = $null; C = $null; D = $null; E = "1.94%" }
    @{ Row = 49; B = $null; C = $null; D = $null; E = $null }
    @{ Row = 50; B = $null; C = $null; D = "0.00002101"; E = "-0.06%" }
    @{ Row = 51; B = $null; C = $null; D = "0.0002001"; E = "-0.06%" }
)

foreach ($r in $rows) {
    if ($r.B -ne $null) { $ws.Cells.Item($r.Row, 2).Value = $r.B }
    if ($r.C -ne $null) { $ws.Cells.Item($r.Row, 3).Value = $r.C }
    if ($r.D -ne $null) {
        $cell = $ws.Cells.Item($r.Row, 4)
        $cell.NumberFormat = "@"
        $cell.Value = $r.D
    }
    if ($r.E -ne $null) {
        $cell = $ws.Cells.Item($r.Row, 5)
        $cell.NumberFormat = "@"
        $cell.Value = $r.E
    }
}
